# ADD results from server
# Update computed result values on the "2025", "2030" and "2035" sheets
# (row 2 = the single data row) to reflect the latest server run.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.005765262920410562
$ws.Range("E2").Value = 0.3690121261547427
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3948917987044045
$ws.Range("L2").Value = 0.570238
$ws.Range("M2").Value = 0.08267716666666668
$ws.Range("N2").Value = 12.85277051881641
$ws.Range("O2").Value = 3.502784081771334

# --- Sheet "2030" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.01740867925716962
$ws.Range("B2").Value = 0.05207980094209073
$ws.Range("E2").Value = 0.2116798874185796
$ws.Range("I2").Value = 0.5409518166666667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04749408333333331
$ws.Range("N2").Value = 5.037181453793078
$ws.Range("O2").Value = 2.225940329859119

# --- Sheet "2035" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.09007993945815861
$ws.Range("B2").Value = 0.02341521357332227
$ws.Range("E2").Value = 0.1667813541380274
$ws.Range("I2").Value = 0.4591872012955962
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.05235633333333339
$ws.Range("N2").Value = 8.135963581684248
$ws.Range("O2").Value = 5.185675745307608
